$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 20834318
$ws.Range("J19").Value = 1062.0834
$ws.Range("L19").Value = 1062.0834
$ws.Range("N19").Value = -1412.0834

$ws.Range("H33").Value = 315.70587
$ws.Range("I33").Value = 260.91666
$ws.Range("J33").Value = 447.2
$ws.Range("K33").Value = 260.91666
$ws.Range("L33").Value = 447.2
$ws.Range("M33").Value = -31.91665999999998
$ws.Range("N33").Value = -905.2

$ws.Range("H40").Value = 1401.9
$ws.Range("I40").Value = 1422.2941
$ws.Range("J40").Value = 1286.3334
$ws.Range("K40").Value = 1422.2941
$ws.Range("L40").Value = 1286.3334
$ws.Range("M40").Value = -1247.2941
$ws.Range("N40").Value = -1636.3334

$ws.Range("H53").Value = 138.6923
$ws.Range("I53").Value = 81.181816
$ws.Range("J53").Value = 180.86667
$ws.Range("K53").Value = 81.181816
$ws.Range("L53").Value = 180.86667
$ws.Range("M53").Value = 555.818184
$ws.Range("N53").Value = -1454.86667

$ws.Range("H63").Value = 38271
$ws.Range("J63").Value = 38271
$ws.Range("L63").Value = 38271
$ws.Range("N63").Value = -39519

$ws.Range("H66").Value = 38271
$ws.Range("J66").Value = 38271
$ws.Range("L66").Value = 114813
$ws.Range("N66").Value = -121053

$ws.Range("H112").Value = 1907.15
$ws.Range("I112").Value = 490
$ws.Range("J112").Value = 1981.7368
$ws.Range("K112").Value = 1470
$ws.Range("L112").Value = 5945.2104
$ws.Range("M112").Value = -362
$ws.Range("N112").Value = -8161.2104

$ws.Range("H132").Value = 1894
$ws.Range("I132").Value = 1938.4138
$ws.Range("K132").Value = 5815.2414
$ws.Range("M132").Value = -3285.2414

$ws.Range("H135").Value = 88236770
$ws.Range("I135").Value = 35715944
$ws.Range("J135").Value = 333333920
$ws.Range("K135").Value = 321443496
$ws.Range("L135").Value = 3000005280
$ws.Range("M135").Value = -321440961
$ws.Range("N135").Value = -3000010350

$ws.Range("H137").Value = 2096.5833
$ws.Range("I137").Value = 1659.7858
$ws.Range("J137").Value = 3625.375
$ws.Range("K137").Value = 4979.357400000001
$ws.Range("L137").Value = 10876.125
$ws.Range("M137").Value = -2429.357400000001
$ws.Range("N137").Value = -15976.125

$ws.Range("H138").Value = 3132.1904
$ws.Range("I138").Value = 1408.1818
$ws.Range("J138").Value = 5028.6
$ws.Range("K138").Value = 4224.5454
$ws.Range("L138").Value = 15085.8
$ws.Range("M138").Value = 915.4546
$ws.Range("N138").Value = -25365.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 11014
$ws.Range("I26").Value = 11014
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 11014
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -10722

$ws.Range("H80").Value = 208.85
$ws.Range("J80").Value = 214.57895
$ws.Range("L80").Value = 214.57895
$ws.Range("N80").Value = -2210.57895

$ws.Range("H83").Value = 208.85
$ws.Range("J83").Value = 214.57895
$ws.Range("L83").Value = 1072.89475
$ws.Range("N83").Value = -11056.89475

$ws.Range("H96").Value = 16627
$ws.Range("I96").Value = 8304
$ws.Range("J96").Value = 24950
$ws.Range("K96").Value = 8304
$ws.Range("L96").Value = 24950
$ws.Range("M96").Value = -5558
$ws.Range("N96").Value = -30442

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2544.0244
$ws.Range("I31").Value = 1826.2759
$ws.Range("J31").Value = 4278.5835
$ws.Range("K31").Value = 1826.2759
$ws.Range("L31").Value = 4278.5835
$ws.Range("M31").Value = -1531.2759
$ws.Range("N31").Value = -4868.5835

$ws.Range("H34").Value = 2544.0244
$ws.Range("I34").Value = 1826.2759
$ws.Range("J34").Value = 4278.5835
$ws.Range("K34").Value = 1826.2759
$ws.Range("L34").Value = 4278.5835
$ws.Range("M34").Value = -1624.2759
$ws.Range("N34").Value = -4682.5835

$ws.Range("H76").Value = 9615.385
$ws.Range("I76").Value = 9615.385
$ws.Range("K76").Value = 9615.385
$ws.Range("M76").Value = -9300.385

$ws.Range("H79").Value = 9615.385
$ws.Range("I79").Value = 9615.385
$ws.Range("K79").Value = 9615.385
$ws.Range("M79").Value = -8523.385

$ws.Range("H134").Value = 3135.628
$ws.Range("I134").Value = 2302.7
$ws.Range("J134").Value = 3859.913
$ws.Range("K134").Value = 6908.099999999999
$ws.Range("L134").Value = 11579.739
$ws.Range("M134").Value = -4373.099999999999
$ws.Range("N134").Value = -16649.739

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 584.5833
$ws.Range("I114").Value = 596.5
$ws.Range("J114").Value = 578.625
$ws.Range("K114").Value = 1789.5
$ws.Range("L114").Value = 1735.875
$ws.Range("M114").Value = 1464.5
$ws.Range("N114").Value = -8243.875

$ws.Range("H117").Value = 1432
$ws.Range("J117").Value = 1782.75
$ws.Range("L117").Value = 5348.25
$ws.Range("N117").Value = -12232.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 10007
$ws.Range("J21").Value = 10007
$ws.Range("L21").Value = 10007
$ws.Range("N21").Value = -10353

$ws.Range("H30").Value = 10007
$ws.Range("J30").Value = 10007
$ws.Range("L30").Value = 10007
$ws.Range("N30").Value = -10217

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = 0

$ws.Range("H7").Value = 5332.095
$ws.Range("I7").Value = 4840.737
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 4840.737
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -4728.737
$ws.Range("N7").Value = -10224

$ws.Range("H16").Value = 1055.1
$ws.Range("I16").Value = 815.7143
$ws.Range("J16").Value = 1613.6666
$ws.Range("K16").Value = 815.7143
$ws.Range("L16").Value = 1613.6666
$ws.Range("M16").Value = -645.7143
$ws.Range("N16").Value = -1953.6666

$ws.Range("H126").Value = 5332.095
$ws.Range("I126").Value = 4840.737
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 14522.211
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -12052.211
$ws.Range("N126").Value = -34940

$ws.Range("H136").Value = 3691.5625
$ws.Range("I136").Value = 2051.244
$ws.Range("J136").Value = 6615.609
$ws.Range("K136").Value = 6153.732
$ws.Range("L136").Value = 19846.827
$ws.Range("M136").Value = -3603.732
$ws.Range("N136").Value = -24946.827

$ws.Range("H140").Value = 52311.938
$ws.Range("J140").Value = 52311.938
$ws.Range("L140").Value = 52311.938
$ws.Range("N140").Value = -62671.938

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 15
$ws.Range("I21").Value = 15
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 15
$ws.Range("L21").ClearContents()
$ws.Range("N21").Value = 0
$ws.Range("M21").Value = 220

$ws.Range("H25").Value = 17400
$ws.Range("I25").Value = 15000
$ws.Range("K25").Value = 15000
$ws.Range("M25").Value = -14707

$ws.Range("H30").Value = 7154.5
$ws.Range("I30").Value = 6272.6665
$ws.Range("K30").Value = 6272.6665
$ws.Range("M30").Value = -6165.6665

$ws.Range("H35").Value = 15
$ws.Range("I35").Value = 15
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 15
$ws.Range("L35").ClearContents()
$ws.Range("N35").Value = 0
$ws.Range("M35").Value = 275

$ws.Range("H40").Value = 32047.555
$ws.Range("I40").Value = 9800
$ws.Range("J40").Value = 38404
$ws.Range("K40").Value = 9800
$ws.Range("L40").Value = 38404
$ws.Range("M40").Value = -9651
$ws.Range("N40").Value = -38702

$ws.Range("H81").Value = 2713.2104
$ws.Range("I81").Value = 1631.375
$ws.Range("K81").Value = 3262.75
$ws.Range("M81").Value = -2201.75

$ws.Range("H84").Value = 2713.2104
$ws.Range("I84").Value = 1631.375
$ws.Range("K84").Value = 16313.75
$ws.Range("M84").Value = -11009.75

$ws.Range("H132").Value = 3249
$ws.Range("I132").Value = 1476
$ws.Range("J132").Value = 5022
$ws.Range("K132").Value = 4428
$ws.Range("L132").Value = 15066
$ws.Range("M132").Value = -1898
$ws.Range("N132").Value = -20126

$ws.Range("H137").Value = 47500
$ws.Range("I137").Value = 30000
$ws.Range("J137").Value = 65000
$ws.Range("K137").Value = 30000
$ws.Range("L137").Value = 65000
$ws.Range("M137").Value = -24900
$ws.Range("N137").Value = -75200
